# Auto-generated edit script: updates market-price derived columns (H-N)
# across the 8 item-category worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to reflect refreshed market data from the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5171
$ws.Range("I40").Value = 4479
$ws.Range("J40").Value = 5555.4443
$ws.Range("K40").Value = 4479
$ws.Range("L40").Value = 5555.4443
$ws.Range("M40").Value = -4304
$ws.Range("N40").Value = -5905.4443
$ws.Range("H98").Value = 1908.1515
$ws.Range("I98").Value = 1949.0312
$ws.Range("J98").Value = 600
$ws.Range("K98").Value = 1949.0312
$ws.Range("L98").Value = 600
$ws.Range("M98").Value = -451.0311999999999
$ws.Range("N98").Value = -3596
$ws.Range("H116").Value = 5261.4
$ws.Range("J116").Value = 7904.5
$ws.Range("L116").Value = 7904.5
$ws.Range("N116").Value = -14788.5
$ws.Range("H122").Value = 1908.1515
$ws.Range("I122").Value = 1949.0312
$ws.Range("J122").Value = 600
$ws.Range("K122").Value = 5847.0936
$ws.Range("L122").Value = 1800
$ws.Range("M122").Value = -3397.0936
$ws.Range("N122").Value = -6700
$ws.Range("H131").Value = 2507021.2
$ws.Range("I131").Value = 3334361.8
$ws.Range("J131").Value = 25000
$ws.Range("K131").Value = 10003085.4
$ws.Range("L131").Value = 75000
$ws.Range("M131").Value = -9998045.399999999
$ws.Range("N131").Value = -85080
$ws.Range("H132").Value = 6543.5386
$ws.Range("I132").Value = 6214.25
$ws.Range("J132").Value = 10495
$ws.Range("K132").Value = 18642.75
$ws.Range("L132").Value = 31485
$ws.Range("M132").Value = -16112.75
$ws.Range("N132").Value = -36545
$ws.Range("H134").Value = 67082.5
$ws.Range("J134").Value = 67082.5
$ws.Range("L134").Value = 67082.5
$ws.Range("N134").Value = -77222.5
$ws.Range("H135").Value = 665.5625
$ws.Range("I135").Value = 576.6
$ws.Range("K135").Value = 5189.400000000001
$ws.Range("M135").Value = -2654.400000000001
$ws.Range("H137").Value = 5850.2593
$ws.Range("I137").Value = 2666.2222
$ws.Range("J137").Value = 7442.278
$ws.Range("K137").Value = 7998.6666
$ws.Range("L137").Value = 22326.834
$ws.Range("M137").Value = -5448.6666
$ws.Range("N137").Value = -27426.834
$ws.Range("H138").Value = 4613.2935
$ws.Range("I138").Value = 2321.5715
$ws.Range("J138").Value = 5504.5186
$ws.Range("K138").Value = 6964.7145
$ws.Range("L138").Value = 16513.5558
$ws.Range("M138").Value = -1824.7145
$ws.Range("N138").Value = -26793.5558
$ws.Range("H141").Value = 2253.4546
$ws.Range("I141").Value = 2253.4546
$ws.Range("K141").Value = 6760.3638
$ws.Range("M141").Value = -1580.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1145410.9
$ws.Range("I32").Value = 530520.3
$ws.Range("K32").Value = 530520.3
$ws.Range("M32").Value = -530233.3
$ws.Range("H45").Value = 100164880
$ws.Range("I45").Value = 206105
$ws.Range("K45").Value = 206105
$ws.Range("M45").Value = -205728
$ws.Range("H61").Value = 3237.6667
$ws.Range("I61").Value = 1224
$ws.Range("K61").Value = 1224
$ws.Range("M61").Value = -1012
$ws.Range("H74").Value = 66677480
$ws.Range("I74").Value = 6597.6
$ws.Range("J74").Value = 100012920
$ws.Range("K74").Value = 6597.6
$ws.Range("L74").Value = 100012920
$ws.Range("M74").Value = -5723.6
$ws.Range("N74").Value = -100014668
$ws.Range("H77").Value = 66677480
$ws.Range("I77").Value = 6597.6
$ws.Range("J77").Value = 100012920
$ws.Range("K77").Value = 32988
$ws.Range("L77").Value = 500064600
$ws.Range("M77").Value = -28620
$ws.Range("N77").Value = -500073336
$ws.Range("H122").Value = 4237.8823
$ws.Range("I122").Value = 3822.182
$ws.Range("K122").Value = 11466.546
$ws.Range("M122").Value = -9016.545999999998
$ws.Range("H132").Value = 8550508
$ws.Range("J132").Value = 3500
$ws.Range("L132").Value = 10500
$ws.Range("N132").Value = -15560
$ws.Range("H136").Value = 3237.6667
$ws.Range("I136").Value = 1224
$ws.Range("K136").Value = 3672
$ws.Range("M136").Value = -1122

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3835.5
$ws.Range("I134").Value = 4002.4
$ws.Range("K134").Value = 12007.2
$ws.Range("M134").Value = -9472.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1584831.1
$ws.Range("I31").Value = 1214.0625
$ws.Range("J31").Value = 1987019.6
$ws.Range("K31").Value = 1214.0625
$ws.Range("L31").Value = 1987019.6
$ws.Range("M31").Value = -919.0625
$ws.Range("N31").Value = -1987609.6
$ws.Range("H34").Value = 1584831.1
$ws.Range("I34").Value = 1214.0625
$ws.Range("J34").Value = 1987019.6
$ws.Range("K34").Value = 1214.0625
$ws.Range("L34").Value = 1987019.6
$ws.Range("M34").Value = -1012.0625
$ws.Range("N34").Value = -1987423.6
$ws.Range("H99").Value = 100003740
$ws.Range("I99").Value = 166669800
$ws.Range("K99").Value = 166669800
$ws.Range("M99").Value = -166668302
$ws.Range("H107").Value = 19231312
$ws.Range("I107").Value = 31250294
$ws.Range("K107").Value = 31250294
$ws.Range("M107").Value = -31248374
$ws.Range("H126").Value = 100003740
$ws.Range("I126").Value = 166669800
$ws.Range("K126").Value = 500009400
$ws.Range("M126").Value = -500006930
$ws.Range("H132").Value = 2866.5757
$ws.Range("I132").Value = 2440.9092
$ws.Range("K132").Value = 7322.7276
$ws.Range("M132").Value = -4792.7276
$ws.Range("H134").Value = 3113.9167
$ws.Range("I134").Value = 3205.24
$ws.Range("J134").Value = 2906.3635
$ws.Range("K134").Value = 9615.719999999999
$ws.Range("L134").Value = 8719.0905
$ws.Range("M134").Value = -7080.719999999999
$ws.Range("N134").Value = -13789.0905

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3851210.8
$ws.Range("J68").Value = 7150989
$ws.Range("L68").Value = 21452967
$ws.Range("N68").Value = -21454589
$ws.Range("H71").Value = 3851210.8
$ws.Range("J71").Value = 7150989
$ws.Range("L71").Value = 64358901
$ws.Range("N71").Value = -64367013
$ws.Range("H107").Value = 5818.129
$ws.Range("I107").Value = 836.5
$ws.Range("J107").Value = 6556.148
$ws.Range("K107").Value = 2509.5
$ws.Range("L107").Value = 19668.444
$ws.Range("M107").Value = -589.5
$ws.Range("N107").Value = -23508.444
$ws.Range("H131").Value = 1791000.6
$ws.Range("J131").Value = 5884588
$ws.Range("L131").Value = 17653764
$ws.Range("N131").Value = -17663844
$ws.Range("H132").Value = 9424.954
$ws.Range("J132").Value = 9335.25
$ws.Range("L132").Value = 84017.25
$ws.Range("N132").Value = -89077.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 245000000
$ws.Range("I3").Value = 245000000
$ws.Range("K3").Value = 245000000
$ws.Range("M3").Value = -244999884
$ws.Range("H10").Value = 100000000
$ws.Range("I10").Value = 100000000
$ws.Range("K10").Value = 100000000
$ws.Range("M10").Value = -99999831
$ws.Range("H102").Value = 125013340
$ws.Range("I102").Value = 142872220
$ws.Range("K102").Value = 142872220
$ws.Range("M102").Value = -142870598
$ws.Range("H132").Value = 3014.6667
$ws.Range("I132").Value = 3000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4117.4
$ws.Range("I7").Value = 4396.75
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 4396.75
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -4284.75
$ws.Range("N7").Value = -3224
$ws.Range("H16").Value = 1503.5
$ws.Range("I16").Value = 1509.2
$ws.Range("K16").Value = 1509.2
$ws.Range("M16").Value = -1339.2
$ws.Range("H40").Value = 70749.664
$ws.Range("I40").Value = 102749.5
$ws.Range("J40").Value = 6750
$ws.Range("K40").Value = 102749.5
$ws.Range("L40").Value = 6750
$ws.Range("M40").Value = -102613.5
$ws.Range("N40").Value = -7022
$ws.Range("H46").Value = 5340.8335
$ws.Range("I46").Value = 4898.5
$ws.Range("J46").Value = 5562
$ws.Range("K46").Value = 4898.5
$ws.Range("L46").Value = 5562
$ws.Range("M46").Value = -4710.5
$ws.Range("N46").Value = -5938
$ws.Range("H55").Value = 584.4783
$ws.Range("I55").Value = 536.3333
$ws.Range("K55").Value = 536.3333
$ws.Range("M55").Value = -363.3333
$ws.Range("H100").Value = 2599.6
$ws.Range("I100").Value = 2249.75
$ws.Range("K100").Value = 2249.75
$ws.Range("M100").Value = -1708.75
$ws.Range("H122").Value = 6007.04
$ws.Range("I122").Value = 5884.357
$ws.Range("K122").Value = 17653.071
$ws.Range("M122").Value = -15203.071
$ws.Range("H126").Value = 4117.4
$ws.Range("I126").Value = 4396.75
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 13190.25
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -10720.25
$ws.Range("N126").Value = -13940
$ws.Range("H132").Value = 10637.654
$ws.Range("I132").Value = 9579.950000000001
$ws.Range("K132").Value = 28739.85
$ws.Range("M132").Value = -26209.85
$ws.Range("H136").Value = 6836.3477
$ws.Range("I136").Value = 5020.4165
$ws.Range("K136").Value = 15061.2495
$ws.Range("M136").Value = -12511.2495

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 773.75
$ws.Range("I96").Value = 773.75
$ws.Range("K96").Value = 773.75
$ws.Range("M96").Value = 599.25
$ws.Range("H122").Value = 14707350
$ws.Range("I122").Value = 1424.5
$ws.Range("K122").Value = 4273.5
$ws.Range("M122").Value = -1823.5
$ws.Range("H132").Value = 5633.032
$ws.Range("I132").Value = 5462.5
$ws.Range("J132").Value = 5943.091
$ws.Range("K132").Value = 16387.5
$ws.Range("L132").Value = 17829.273
$ws.Range("M132").Value = -13857.5
$ws.Range("N132").Value = -22889.273
